# Remove the trailing bullet paragraph "Replaces HTML and CSS with JSX"
# from the body placeholder on slide 2 of the presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the shape whose text contains the bullet we need to remove,
# rather than hard-coding a shape index.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -like "*Replaces HTML and CSS with JSX*") {
            $targetShape = $sh
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# The paragraph we want gone is the very last paragraph in the text body.
$paraCount = $tr.Paragraphs().Count
$lastPara = $tr.Paragraphs($paraCount, 1)

# Select the last paragraph's own text plus one extra character past its end.
# Because this paragraph is the final paragraph in the story it carries no
# trailing paragraph mark of its own; reaching one character beyond its text
# pulls in the paragraph break that precedes it so the whole paragraph node
# (not just its run text) is deleted and the previous paragraph is left
# completely untouched.
$delStart = $lastPara.Start
$delLen = $lastPara.Length + 1
$tr.Characters($delStart, $delLen).Delete()
